$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (riparian) values:
# B7: 0 -> 120 (kept as text, matching the existing text-number formatting in column B)
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "120"

# D7: 0 -> 0.94 (numeric)
$ws.Range("D7").Value = 0.94
